# ============================================================================
# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
#  - new sheet "Player Info"  (inserted before "ODI Batting")
#  - "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (header + values trimmed to
#    the numeric match code), and the stray blank INNING_NUMBER cells on the
#    "did not bat" rows are cleared out entirely
#  - "ODI Bowling": same MATCH_CARD_LINK -> MATCH_CODE treatment
#  - new sheet "ODI Batting Extra" (inserted after "ODI Bowling") holding the
#    newly scraped per-match batting detail
# ============================================================================

$wb = $excel.ActiveWorkbook

function Set-HeaderRow($ws, $row, $headers) {
    $col = 1
    foreach ($h in $headers) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $h
        $col = $col + 1
    }
    $count = $headers.Length
    $headerRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, $count))
    $headerRange.Font.Bold = $true
    $headerRange.HorizontalAlignment = -4108   # xlCenter
    $headerRange.VerticalAlignment = -4160     # xlTop
    $headerRange.Borders.LineStyle = 1
}

# Pull the numeric MatchCode out of a howstat scorecard URL, e.g.
# ".../MatchScorecard_ODI.asp?MatchCode=2987" -> "2987"
function Get-MatchCodeFromLink($link) {
    $parts = $link.Split("=")
    return $parts[$parts.Length - 1]
}

# ----------------------------------------------------------------------------
# 1) New "Player Info" sheet, inserted before "ODI Batting"
# ----------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

Set-HeaderRow $playerInfo 1 @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "3718"
$playerInfo.Cells.Item(2, 2).Value = "Devon Cuthbert Thomas"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium"

# ----------------------------------------------------------------------------
# 2) "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE
# ----------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$lastRow = $battingSheet.UsedRange.Rows.Count
$battingSheet.Range($battingSheet.Cells.Item(2, 4), $battingSheet.Cells.Item($lastRow, 4)).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $link = $battingSheet.Cells.Item($r, 4).Text
    $code = Get-MatchCodeFromLink $link
    $battingSheet.Cells.Item($r, 4).Value = $code
}

# rows where the player did not bat carry a stray blank INNING_NUMBER cell
# that the scrape now drops entirely
$battingSheet.Cells.Item(2, 2).ClearContents()
$battingSheet.Cells.Item(6, 2).ClearContents()

# ----------------------------------------------------------------------------
# 3) "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE
# ----------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$lastRowBowl = $bowlingSheet.UsedRange.Rows.Count
$bowlingSheet.Range($bowlingSheet.Cells.Item(2, 2), $bowlingSheet.Cells.Item($lastRowBowl, 2)).NumberFormat = "@"

for ($r = 2; $r -le $lastRowBowl; $r++) {
    $link = $bowlingSheet.Cells.Item($r, 2).Text
    $code = Get-MatchCodeFromLink $link
    $bowlingSheet.Cells.Item($r, 2).Value = $code
}

# ----------------------------------------------------------------------------
# 4) New "ODI Batting Extra" sheet, inserted after "ODI Bowling"
# ----------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowlingSheet)
$extra.Name = "ODI Batting Extra"

Set-HeaderRow $extra 1 @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("2989", 8, "2", "0", "11.69%", "NO"),
    @("3241", 7, "1", "0", "6.76%",  "NO"),
    @("3247", $null, $null, $null, $null, "NO"),
    @("3253", 7, $null, $null, $null, "NO"),
    @("3261", 7, "0", "0", "0.73%",  "NO"),
    @("3270", 5, "1", "0", "4.00%",  "NO"),
    @("3276", 6, "0", "0", "1.06%",  "NO"),
    @("3277", 7, "0", "0", $null,    "NO"),
    @("3429", 7, "3", "1", "14.02%", "NO"),
    @("3430", 7, "1", "0", "8.30%",  "NO"),
    @("3450", 6, "1", "0", "8.04%",  "NO"),
    @("3451", 7, "0", "0", $null,    "NO"),
    @("3452", 7, "1", "0", "5.26%",  "NO"),
    @("3453", 7, "1", "0", "4.74%",  "NO"),
    @("3454", $null, $null, $null, $null, "NO"),
    @("3471", 7, "0", "0", "4.29%",  "NO"),
    @("3472", 7, "0", "0", $null,    "NO"),
    @("3473", 2, "3", "0", "6.55%",  "NO"),
    @("3474", 7, "1", "0", "3.18%",  "NO"),
    @("3475", 7, "1", "1", "7.39%",  "NO")
)

$rowCount = $extraRows.Length
$extra.Range($extra.Cells.Item(2, 1), $extra.Cells.Item($rowCount + 1, 1)).NumberFormat = "@"
$extra.Range($extra.Cells.Item(2, 3), $extra.Cells.Item($rowCount + 1, 4)).NumberFormat = "@"
$extra.Range($extra.Cells.Item(2, 5), $extra.Cells.Item($rowCount + 1, 5)).NumberFormat = "@"

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($row[2] -ne $null) {
        $extra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -ne $null) {
        $extra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($row[4] -ne $null) {
        $extra.Cells.Item($r, 5).Value = $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

Write-Host "Edit complete. Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host " - $($s.Name)"
}
